# GINF2 Modules update
# Fills in the module/teacher/component table (rows 2-13) and renames the
# header cells in row 1, then tidies up the column widths and selection to
# match the author's final layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Intitulé"
$ws.Range("C1").Value = "Chef Module"
$ws.Range("D1").Value = "Composants"

# --- Data rows --------------------------------------------------------------
$ws.Range("A2").Value = "GINF31"
$ws.Range("B2").Value = "POO & XML"
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Java,XML"

$ws.Range("A3").Value = "GINF32"
$ws.Range("B3").Value = "Qualité et approche processus"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "Qualité, Cycle de vie Logiciel, Optimisation des processus"

$ws.Range("A4").Value = "GINF33"
$ws.Range("B4").Value = "Modélisation OO & IHM"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "UML,IHM"

$ws.Range("A5").Value = "GINF34"
$ws.Range("B5").Value = "BD Avancées I"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Optimisation BD, Admin. BD, BD Distribuées"

$ws.Range("A6").Value = "GINF35"
$ws.Range("B6").Value = "Admin. Et prog. Système"
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "Ad. Sys, Prog. Sys"

$ws.Range("A7").Value = "GINF36"
$ws.Range("B7").Value = "Langues et Communication"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "Anglais, Espagnol"

$ws.Range("A8").Value = "GINF41"
$ws.Range("B8").Value = "Technologies distribués"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "J2EE, C#"

$ws.Range("A9").Value = "GINF42"
$ws.Range("B9").Value = "BD Avacncées II & Cloud"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "Gestion des données Complexes, NOSQL, Cloud Computing"

$ws.Range("A10").Value = "GINF43"
$ws.Range("B10").Value = "Traitement de l'image "
$ws.Range("C10").Value = "Badir"
$ws.Range("D10").Value = "Traitement d'image, vision numérique, Processus stochastique"

$ws.Range("A11").Value = "GINF44"
$ws.Range("B11").Value = "Prog. Déclarative et TAV"
$ws.Range("C11").Value = "Ezzine"
$ws.Range("D11").Value = "Prog. Déclarative, Technique algorithmique avancée"

$ws.Range("A12").Value = "GINF45"
$ws.Range("B12").Value = "Sécurité des systèmes & Crypto."
$ws.Range("C12").Value = "Ben Achrab"
$ws.Range("D12").Value = "Sécurité des systèmes, Cryptographie"

$ws.Range("A13").Value = "GINF46"
$ws.Range("B13").Value = "Management de l'entreprise"
$ws.Range("C13").Value = "EL Haddad"
$ws.Range("D13").Value = "Economie & Compta 2, Projets collectifs, Management de projet"

# --- Column widths (matches the widened Intitulé / Chef Module / Composants columns) ---
$ws.Columns.Item(2).ColumnWidth = 34.666666666666664
$ws.Columns.Item(3).ColumnWidth = 20.166666666666668
$ws.Columns.Item(4).ColumnWidth = 49.666666666666664

# --- Final selection --------------------------------------------------------
$ws.Range("C4").Select()
